# Auto-generated edit script applying the Kujata_Profits market-data refresh
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets (H/I/J/K/L/M/N price & profit columns).
$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2530
$ws.Range("J18").Value = 1900
$ws.Range("L18").Value = 1900
$ws.Range("N18").Value = -2468

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 5180
$ws.Range("I111").Value = 5300
$ws.Range("K111").Value = 15900
$ws.Range("M111").Value = -12833

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 877.0540999999999
$ws.Range("I129").Value = 481.77777
$ws.Range("J129").Value = 1004.1071
$ws.Range("K129").Value = 1445.33331
$ws.Range("L129").Value = 3012.3213
$ws.Range("M129").Value = 3554.66669
$ws.Range("N129").Value = -13012.3213

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7755672.5
$ws.Range("I132").Value = 8774017
$ws.Range("J132").Value = 16258
$ws.Range("K132").Value = 26322051
$ws.Range("L132").Value = 48774
$ws.Range("M132").Value = -26319521
$ws.Range("N132").Value = -53834

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1490.0588
$ws.Range("I137").Value = 1431.3
$ws.Range("K137").Value = 4293.9
$ws.Range("M137").Value = -1743.9

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1234
$ws.Range("I138").Value = 812.68085
$ws.Range("J138").Value = 1716.9756
$ws.Range("K138").Value = 2438.04255
$ws.Range("L138").Value = 5150.9268
$ws.Range("M138").Value = 2701.95745
$ws.Range("N138").Value = -15430.9268

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 893.2222
$ws.Range("I141").Value = 893.2222
$ws.Range("K141").Value = 2679.6666
$ws.Range("M141").Value = 2500.3334

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3251.4375
$ws.Range("I32").Value = 2845.689
$ws.Range("K32").Value = 2845.689
$ws.Range("M32").Value = -2558.689

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1321
$ws.Range("I61").Value = 948.1667
$ws.Range("J61").Value = 2066.6667
$ws.Range("K61").Value = 948.1667
$ws.Range("L61").Value = 2066.6667
$ws.Range("M61").Value = -736.1667
$ws.Range("N61").Value = -2490.6667

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1192.091
$ws.Range("I74").Value = 776.5
$ws.Range("J74").Value = 1690.8
$ws.Range("K74").Value = 776.5
$ws.Range("L74").Value = 1690.8
$ws.Range("M74").Value = 97.5
$ws.Range("N74").Value = -3438.8

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1192.091
$ws.Range("I77").Value = 776.5
$ws.Range("J77").Value = 1690.8
$ws.Range("K77").Value = 3882.5
$ws.Range("L77").Value = 8454
$ws.Range("M77").Value = 485.5
$ws.Range("N77").Value = -17190

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2354.2
$ws.Range("I132").Value = 1970.625
$ws.Range("K132").Value = 5911.875
$ws.Range("M132").Value = -3381.875

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1321
$ws.Range("I136").Value = 948.1667
$ws.Range("J136").Value = 2066.6667
$ws.Range("K136").Value = 2844.5001
$ws.Range("L136").Value = 6200.000100000001
$ws.Range("M136").Value = -294.5001000000002
$ws.Range("N136").Value = -11300.0001

# ARM row 138
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 64440
$ws.Range("J138").Value = 64440
$ws.Range("L138").Value = 64440
$ws.Range("N138").Value = -74720

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1433.3334
$ws.Range("I107").Value = 1075
$ws.Range("K107").Value = 1075
$ws.Range("M107").Value = 845

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7524.5713
$ws.Range("I134").Value = 947.5
$ws.Range("J134").Value = 13503.728
$ws.Range("K134").Value = 2842.5
$ws.Range("L134").Value = 40511.18399999999
$ws.Range("M134").Value = -307.5
$ws.Range("N134").Value = -45581.18399999999

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 371.14285
$ws.Range("I22").Value = 275
$ws.Range("J22").Value = 499.33334
$ws.Range("K22").Value = 275
$ws.Range("L22").Value = 499.33334
$ws.Range("M22").Value = 75
$ws.Range("N22").Value = -1199.33334

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1735
$ws.Range("I58").Value = 1392.3636
$ws.Range("J58").Value = 2273.4285
$ws.Range("K58").Value = 1392.3636
$ws.Range("L58").Value = 2273.4285
$ws.Range("M58").Value = -1189.3636
$ws.Range("N58").Value = -2679.4285

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5146959.5
$ws.Range("I86").Value = 13336572
$ws.Range("K86").Value = 13336572
$ws.Range("M86").Value = -13335449

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 5146959.5
$ws.Range("I89").Value = 13336572
$ws.Range("K89").Value = 66682860
$ws.Range("M89").Value = -66677244

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 15752.625
$ws.Range("I132").Value = 27005.75
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 81017.25
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -78487.25
$ws.Range("N132").Value = -18558.5

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1958.4166
$ws.Range("I134").Value = 1997.7778
$ws.Range("K134").Value = 5993.3334
$ws.Range("M134").Value = -3458.3334

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1735
$ws.Range("I136").Value = 1392.3636
$ws.Range("J136").Value = 2273.4285
$ws.Range("K136").Value = 4177.0908
$ws.Range("L136").Value = 6820.2855
$ws.Range("M136").Value = -1627.0908
$ws.Range("N136").Value = -11920.2855

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 558.087
$ws.Range("I107").Value = 709.9167
$ws.Range("J107").Value = 392.45456
$ws.Range("K107").Value = 709.9167
$ws.Range("L107").Value = 392.45456
$ws.Range("M107").Value = 1210.0833
$ws.Range("N107").Value = -4232.45456

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1631.2222
$ws.Range("J22").Value = 1760.125
$ws.Range("L22").Value = 1760.125
$ws.Range("N22").Value = -2350.125

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1631.2222
$ws.Range("J27").Value = 1760.125
$ws.Range("L27").Value = 1760.125
$ws.Range("N27").Value = -1974.125

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5277.6665
$ws.Range("I136").Value = 7046.1177
$ws.Range("J136").Value = 2271.3
$ws.Range("K136").Value = 21138.3531
$ws.Range("L136").Value = 6813.900000000001
$ws.Range("M136").Value = -18588.3531
$ws.Range("N136").Value = -11913.9

# WVR row 120
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 32500
$ws.Range("J120").Value = 32500
$ws.Range("L120").Value = 32500
$ws.Range("N120").Value = -42176

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2241.3823
$ws.Range("I132").Value = 1696.8518
$ws.Range("J132").Value = 4341.7144
$ws.Range("K132").Value = 5090.555399999999
$ws.Range("L132").Value = 13025.1432
$ws.Range("M132").Value = -2560.555399999999
$ws.Range("N132").Value = -18085.1432

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 55357.5
$ws.Range("J135").Value = 55357.5
$ws.Range("L135").Value = 55357.5
$ws.Range("N135").Value = -65497.5
